$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 with the new Russian translation text
$ws.Range("B1").Value = "9.5.2 Количество исследователей (в эквиваленте полной занятости) на миллион жителей"

# Add new column Q for year 2023 with value 631, matching the formatting of column P
$ws.Range("P4:P5").Copy() | Out-Null
$ws.Range("Q4:Q5").PasteSpecial(-4122) | Out-Null

$ws.Range("Q4").Value = 2023
$ws.Range("Q5").Value = 631

# Move selection away from P10 (diff removes the explicit <selection> element)
$ws.Range("A1").Select() | Out-Null
